# 2019.10.08. 김동욱 PSP 수정
#
# Fills in the five missing PSP time-log entries for 김동욱 (rows 6-10 of
# his sheet) and leaves the workbook's UI state (active sheet / selected
# cells) the way a person would after typing them in: the "김동욱" tab
# ends up active, with a couple of other tabs showing the cell the editor
# last clicked on while passing through them.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) 김지환 sheet: briefly visited, cursor left on F10.
# ---------------------------------------------------------------------
$wsJihwan = $wb.Worksheets.Item("김지환")
$wsJihwan.Activate()
$wsJihwan.Range("F10").Select()

# ---------------------------------------------------------------------
# 2) 김동욱 sheet: add the five new PSP log rows.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("김동욱")
$ws.Activate()

# Row 6 already carries the date/time/activity formatting used as the
# template for the new rows below - copy it over E:F for rows 7-10 first
# (columns A-D already share the same style across rows 6-10).
$ws.Range("E6:F6").Copy()
$ws.Range("E7:F10").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 6 - 2019.09.25, 19:00-20:00, 0 interruption, 60 min
$ws.Range("A6").Value = 43733
$ws.Range("B6").Value = 0.79166666666666663
$ws.Range("C6").Value = 0.83333333333333337
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 60
$ws.Range("F6").Value = "usecase diagram 작성"

# Row 7 - 2019.09.27, 15:30-17:20, 0 interruption, 110 min
$ws.Range("A7").Value = 43735
$ws.Range("B7").Value = 0.64583333333333337
$ws.Range("C7").Value = 0.72222222222222221
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 110
$ws.Range("F7").Value = "usecase 및 prototype 작성 관련 조모임"

# Row 8 - 2019.10.04, 16:30-17:30, 0 interruption, 60 min
$ws.Range("A8").Value = 43742
$ws.Range("B8").Value = 0.6875
$ws.Range("C8").Value = 0.72916666666666663
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 60
$ws.Range("F8").Value = "spec 작성 관련 조모임"

# Row 9 - 2019.10.05, 17:00-18:00, 0 interruption, 60 min
$ws.Range("A9").Value = 43743
$ws.Range("B9").Value = 0.70833333333333337
$ws.Range("C9").Value = 0.75
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 60
$ws.Range("F9").Value = "추천 시스템 탐색 및 스터디"

# Row 10 - 2019.10.06, 18:00-22:00, 0 interruption, 240 min
$ws.Range("A10").Value = 43744
$ws.Range("B10").Value = 0.75
$ws.Range("C10").Value = 0.91666666666666663
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 240
$ws.Range("F10").Value = "usecase outline 수정 및 specification 작성"

# Leave 김동욱 as the active tab with the cursor on H8, matching where the
# editor ended up after typing in the log.
$ws.Range("H8").Select()
